$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flight Mission Cycle")
$ws.Range("B1").Value = "Duration"
$ws.Activate()
[void]$ws.Range("B1").Select()
